# "New command - Rotate2PL"
# Insert a new alphabetically-ordered command row ("ROTATE2PL") into the
# RFL Tools command listing, between "ROTATE2ALIGN" (row 126) and
# "ROUNDABOUT" (old row 127 / new row 128).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 127, pushing ROUNDABOUT and everything below it
# down by one row.
$ws.Rows("127:127").Insert()

# Fill in the new command name + description.
$ws.Range("A127").Value2 = "ROTATE2PL"
$ws.Range("B127").Value2 = "Rotates blocks to be aligned with selected polylines"

# Match the saved selection/viewport (B128 is now ROUNDABOUT's description cell).
$ws.Range("B128").Select()

Write-Output "Inserted ROTATE2PL at row 127"
